$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: Oyuncu Adı (Player Name)
$ws.Range("A2").Value = "Coby White"
$ws.Range("A3").Value = "LaMelo Ball"
$ws.Range("A4").Value = "Quentin Grimes"
$ws.Range("A5").Value = "Jusuf Nurkic"
$ws.Range("A6").Value = "Robert Williams III"
$ws.Range("A7").Value = "Isaiah Hartenstein"
$ws.Range("A8").Value = "Cade Cunningham"
$ws.Range("A9").Value = "Carlton Carrington"
$ws.Range("A10").Value = "Andrew Wiggins"
$ws.Range("A11").Value = "Ausar Thompson"
$ws.Range("A12").Value = "Precious Achiuwa"
$ws.Range("A13").Value = "Naz Reid"
$ws.Range("A14").Value = "Malik Monk"
$ws.Range("A15").Value = "Derrick White"
$ws.Range("A16").Value = "Onyeka Okongwu"
$ws.Range("A17").Value = "Collin Sexton"
$ws.Range("A18").Value = "Anthony Davis"
$ws.Range("A19").Value = "Damian Lillard"

# Column B: Pozisyon (Position)
$ws.Range("B2").Value = "PG,SG"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("B5").Value = "C"
$ws.Range("B6").Value = "C"
$ws.Range("B7").Value = "C"
$ws.Range("B8").Value = "PG,SG"
$ws.Range("B9").Value = "PG,SG"
$ws.Range("B10").Value = "SF,PF"
$ws.Range("B11").Value = "SF,PF"
$ws.Range("B12").Value = "PF,C"
$ws.Range("B13").Value = "PF,C"
$ws.Range("B14").Value = "PG,SG,SF"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("B16").Value = "PF,C"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("B18").Value = "PF,C"
$ws.Range("B19").Value = "PG"

# Column C: Takım (Team)
$ws.Range("C2").Value = "Chicago Bulls"
$ws.Range("C3").Value = "Charlotte Hornets"
$ws.Range("C4").Value = "Philadelphia 76ers"
$ws.Range("C5").Value = "Charlotte Hornets"
$ws.Range("C6").Value = "Portland Trail Blazers"
$ws.Range("C7").Value = "Oklahoma City Thunder"
$ws.Range("C8").Value = "Detroit Pistons"
$ws.Range("C9").Value = "Washington Wizards"
$ws.Range("C10").Value = "Miami Heat"
$ws.Range("C11").Value = "Detroit Pistons"
$ws.Range("C12").Value = "New York Knicks"
$ws.Range("C13").Value = "Minnesota Timberwolves"
$ws.Range("C14").Value = "Sacramento Kings"
$ws.Range("C15").Value = "Boston Celtics"
$ws.Range("C16").Value = "Atlanta Hawks"
$ws.Range("C17").Value = "Utah Jazz"
$ws.Range("C18").Value = "Dallas Mavericks"
$ws.Range("C19").Value = "Milwaukee Bucks"
